# feat: add 2022-Q4 data
#
# 1. Insert a new "2022-Q4" worksheet (copied from "2022-Q3" so it keeps the
#    same column layout/styles) positioned right after "总计" and before
#    "2022-Q3", populated with the new quarter's fund data.
# 2. Insert a new row into "总计" for the "2022-Q4" summary figures, pushing
#    the existing quarters down by one row.
# 3. Fix the "2021-Q3" sheet's header cell D1 ("基金金额" -> "基金规模") to
#    match the other quarterly sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create the "2022-Q4" sheet from a copy of "2022-Q3"
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# Columns B-G hold text values (fund code/name/ratios) even when they look
# numeric (e.g. "000906", "20.45") - force text formatting before assigning
# so leading zeros / exact decimal text are preserved.
$q4.Range("B2:G3").NumberFormat = "@"

$q4.Range("B2").Value = "270023"
$q4.Range("C2").Value = "广发全球精选股票（QDII）"
$q4.Range("D2").Value = "20.45"
$q4.Range("E2").Value = "82.63"
$q4.Range("F2").Value = "3.59"
$q4.Range("G2").Value = "0.7342"
$q4.Range("H2").Value = 10

$q4.Range("B3").Value = "000906"
$q4.Range("C3").Value = "广发全球精选股票（QDII）美元现汇"
$q4.Range("D3").Value = "20.45"
$q4.Range("E3").Value = "82.63"
$q4.Range("F3").Value = "3.59"
$q4.Range("G3").Value = "0.7342"
$q4.Range("H3").Value = 10

# ---------------------------------------------------------------------------
# Step 2: insert the new "2022-Q4" row into "总计" (row 2), shifting the
# existing quarterly rows down.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# Copy formatting from row 3 (the row that used to be row 2) onto the new
# blank row 2 so styles match the rest of the table.
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 1.47

# ---------------------------------------------------------------------------
# Step 3: relabel "2021-Q3"'s D1 header to match the other quarters.
# ---------------------------------------------------------------------------
$q3_2021 = $wb.Worksheets.Item("2021-Q3")
$q3_2021.Range("D1").Value = "基金规模"
